# Update the acquisition-timestamp column (A) for the existing data rows
# on the "ランサーズ" sheet to reflect the newest scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-07 01:26:36"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
